$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F..V, excluding the already-identical
# G/I goal columns) between row 49 and row 50 -----------------------------
$swapCols = @(6,8,10,11,12,13,14,15,16,17,18,19,20,21,22)
foreach ($c in $swapCols) {
    $v49 = $ws.Cells.Item(49, $c).Value2
    $v50 = $ws.Cells.Item(50, $c).Value2
    $ws.Cells.Item(49, $c).Value = $v50
    $ws.Cells.Item(50, $c).Value = $v49
}

# --- Append the new match as row 57 --------------------------------------
$newRow = 57

# Copy formatting from row 56's styled cells (Indice col A, data_partida col E)
$ws.Cells.Item(56, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

$ws.Cells.Item(56, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = 56
$ws.Cells.Item($newRow, 2).Value = "turkey"
$ws.Cells.Item($newRow, 3).Value = "super-lig"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45194.79166666666
$ws.Cells.Item($newRow, 6).Value = "Hatayspor"
$ws.Cells.Item($newRow, 7).Value = 3
$ws.Cells.Item($newRow, 8).Value = "Trabzonspor"
$ws.Cells.Item($newRow, 9).Value = 2
$ws.Cells.Item($newRow, 10).Value = 3.76
$ws.Cells.Item($newRow, 11).Value = "18/09/2023 18:13"
$ws.Cells.Item($newRow, 12).Value = 3.42
$ws.Cells.Item($newRow, 13).Value = "25/09/2023 18:58"
$ws.Cells.Item($newRow, 14).Value = 3.82
$ws.Cells.Item($newRow, 15).Value = "18/09/2023 18:13"
$ws.Cells.Item($newRow, 16).Value = 3.57
$ws.Cells.Item($newRow, 17).Value = "25/09/2023 18:58"
$ws.Cells.Item($newRow, 18).Value = 1.98
$ws.Cells.Item($newRow, 19).Value = "18/09/2023 18:13"
$ws.Cells.Item($newRow, 20).Value = 2.2
$ws.Cells.Item($newRow, 21).Value = "25/09/2023 18:57"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/turkey/super-lig/hatayspor-trabzonspor/Gb1A9LE0/"
